$d = $word.ActiveDocument

# --- Step 1: locate the "Tela de Cadastros" paragraph and the two break-only
#     paragraphs that immediately precede it (a manual line break paragraph and
#     a page-break paragraph), then collapse all three into a single paragraph
#     that only contains the heading text. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Tela de Cadastros*") {
        $target = $i
        break
    }
}

$headingPara = $d.Paragraphs.Item($target)
$breakPara2 = $d.Paragraphs.Item($target - 1)
$breakPara1 = $d.Paragraphs.Item($target - 2)

# Remove the two break-only paragraphs entirely (this also removes the
# <w:br/> and <w:br w:type="page"/> runs they contained).
$removeRange = $d.Range($breakPara1.Range.Start, $headingPara.Range.Start)
$removeRange.Delete()

# Re-fetch the (now renumbered) heading paragraph and strip the
# <w:lastRenderedPageBreak/> marker by clearing and retyping its text,
# preserving the bold formatting.
$headingPara = $d.Paragraphs.Item($target - 2)
$headingRange = $headingPara.Range
$headingText = $d.Range($headingRange.Start, $headingRange.End - 1)
$headingText.Delete()

$insertPoint = $d.Range($headingPara.Range.Start, $headingPara.Range.Start)
$insertPoint.InsertAfter("Tela de Cadastros")

$newHeadingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$newHeadingRange.Bold = 1

# --- Step 2: move the "_GoBack" bookmark from the end of the following
#     paragraph to the start of the (rebuilt) heading paragraph. ---
$headingPara = $d.Paragraphs.Item($target - 2)
$bmRange = $headingPara.Range
$bmRange.Collapse(1)
$bmRange.Bookmarks.Add("_GoBack")

# --- Step 3: shrink the top/bottom page margins (1417 -> 142 twips, i.e.
#     70.85pt -> 7.1pt). PageSetup margins are expressed in points. ---
$d.Sections.Item(1).PageSetup.TopMargin = 7.1
$d.Sections.Item(1).PageSetup.BottomMargin = 7.1
